$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the quantity value that used to live in B2 (column C now takes over
# the "extra data" role with a Date column instead).
$ws.Range("B2").ClearContents() | Out-Null

# Add the new "Date" header in C1, matching the bold header style already
# used by A1 ("ISBN") and B1 ("Quantity").
$ws.Range("C1").Value = "Date"
$ws.Range("C1").Font.Bold = $true

# Move/collapse the active selection onto the newly added header cell.
$ws.Range("C1").Select() | Out-Null
